# Apply the commit "committing new decision datas":
#  - Append 27 new rows (131-157) of decision data to the "sys2" sheet
#    and the corresponding derived subset of columns to "sys2_test".
#  - Update the saved view state: "sys2_test" becomes the active/selected
#    sheet/tab, with updated selections on both sheets.

$wb  = $excel.ActiveWorkbook
$ws7 = $wb.Worksheets.Item("sys2")
$ws8 = $wb.Worksheets.Item("sys2_test")

# --- New rows 131-157 on sys2 (sheet7) ---
$ws7.Range("A131").Value = 1
$ws7.Range("B131").Value = 1
$ws7.Range("C131").Value = -0.023021396615663402
$ws7.Range("D131").Value = 1
$ws7.Range("E131").Value = 0
$ws7.Range("F131").Value = -0.14308251083017601
$ws7.Range("G131").Value = 0
$ws7.Range("H131").Value = 0
$ws7.Range("I131").Value = 0
$ws7.Range("J131").Value = 1
$ws7.Range("K131").Value = 0
$ws7.Range("L131").Value = 0
$ws7.Range("M131").Value = "right"
$ws7.Range("A132").Value = 0
$ws7.Range("B132").Value = 0
$ws7.Range("C132").Value = 0.060738217814308201
$ws7.Range("D132").Value = 1
$ws7.Range("E132").Value = 1
$ws7.Range("F132").Value = -0.104055189098392
$ws7.Range("G132").Value = 1
$ws7.Range("H132").Value = 0
$ws7.Range("I132").Value = 1
$ws7.Range("J132").Value = 0
$ws7.Range("K132").Value = 1
$ws7.Range("L132").Value = 1
$ws7.Range("M132").Value = "left"
$ws7.Range("A133").Value = 0
$ws7.Range("B133").Value = 1
$ws7.Range("C133").Value = -0.22789823849746699
$ws7.Range("D133").Value = 0
$ws7.Range("E133").Value = 0
$ws7.Range("F133").Value = -0.0085016880706790705
$ws7.Range("G133").Value = 1
$ws7.Range("H133").Value = 1
$ws7.Range("I133").Value = 0
$ws7.Range("J133").Value = 1
$ws7.Range("K133").Value = 0
$ws7.Range("L133").Value = 0
$ws7.Range("M133").Value = "right"
$ws7.Range("A134").Value = 0
$ws7.Range("B134").Value = 0
$ws7.Range("C134").Value = -0.038013183504730999
$ws7.Range("D134").Value = 1
$ws7.Range("E134").Value = 1
$ws7.Range("F134").Value = -0.077462887875195702
$ws7.Range("G134").Value = 1
$ws7.Range("H134").Value = 0
$ws7.Range("I134").Value = 1
$ws7.Range("J134").Value = 0
$ws7.Range("K134").Value = 1
$ws7.Range("L134").Value = 1
$ws7.Range("M134").Value = "left"
$ws7.Range("A135").Value = 1
$ws7.Range("B135").Value = 1
$ws7.Range("C135").Value = -0.053364465595300802
$ws7.Range("D135").Value = 0
$ws7.Range("E135").Value = 0
$ws7.Range("F135").Value = -0.0051054128281936796
$ws7.Range("G135").Value = 0
$ws7.Range("H135").Value = 1
$ws7.Range("I135").Value = 0
$ws7.Range("J135").Value = 1
$ws7.Range("K135").Value = 0
$ws7.Range("L135").Value = 0
$ws7.Range("M135").Value = "right"
$ws7.Range("A136").Value = 0
$ws7.Range("B136").Value = 0
$ws7.Range("C136").Value = 0.060156846920413502
$ws7.Range("D136").Value = 1
$ws7.Range("E136").Value = 1
$ws7.Range("F136").Value = -0.097380944572685701
$ws7.Range("G136").Value = 1
$ws7.Range("H136").Value = 0
$ws7.Range("I136").Value = 1
$ws7.Range("J136").Value = 0
$ws7.Range("K136").Value = 1
$ws7.Range("L136").Value = 1
$ws7.Range("M136").Value = "left"
$ws7.Range("A137").Value = 0
$ws7.Range("B137").Value = 1
$ws7.Range("C137").Value = -0.21314620281322399
$ws7.Range("D137").Value = 0
$ws7.Range("E137").Value = 0
$ws7.Range("F137").Value = 0.050115291735105302
$ws7.Range("G137").Value = 1
$ws7.Range("H137").Value = 1
$ws7.Range("I137").Value = 0
$ws7.Range("J137").Value = 1
$ws7.Range("K137").Value = 0
$ws7.Range("L137").Value = 0
$ws7.Range("M137").Value = "right"
$ws7.Range("A138").Value = 0
$ws7.Range("B138").Value = 0
$ws7.Range("C138").Value = 0.091723723851644304
$ws7.Range("D138").Value = 1
$ws7.Range("E138").Value = 1
$ws7.Range("F138").Value = -0.072693457748126605
$ws7.Range("G138").Value = 1
$ws7.Range("H138").Value = 0
$ws7.Range("I138").Value = 1
$ws7.Range("J138").Value = 0
$ws7.Range("K138").Value = 1
$ws7.Range("L138").Value = 1
$ws7.Range("M138").Value = "left"
$ws7.Range("A139").Value = 1
$ws7.Range("B139").Value = 1
$ws7.Range("C139").Value = -0.032169778198787198
$ws7.Range("D139").Value = 0
$ws7.Range("E139").Value = 0
$ws7.Range("F139").Value = -0.0027072716051931298
$ws7.Range("G139").Value = 0
$ws7.Range("H139").Value = 1
$ws7.Range("I139").Value = 0
$ws7.Range("J139").Value = 1
$ws7.Range("K139").Value = 0
$ws7.Range("L139").Value = 0
$ws7.Range("M139").Value = "right"
$ws7.Range("A140").Value = 0
$ws7.Range("B140").Value = 0
$ws7.Range("C140").Value = 0.061722051235653497
$ws7.Range("D140").Value = 1
$ws7.Range("E140").Value = 1
$ws7.Range("F140").Value = -0.10005958715931799
$ws7.Range("G140").Value = 1
$ws7.Range("H140").Value = 0
$ws7.Range("I140").Value = 1
$ws7.Range("J140").Value = 0
$ws7.Range("K140").Value = 1
$ws7.Range("L140").Value = 1
$ws7.Range("M140").Value = "left"
$ws7.Range("A141").Value = 1
$ws7.Range("B141").Value = 1
$ws7.Range("C141").Value = -0.0266094643094196
$ws7.Range("D141").Value = 1
$ws7.Range("E141").Value = 0
$ws7.Range("F141").Value = -0.14829029022474799
$ws7.Range("G141").Value = 0
$ws7.Range("H141").Value = 0
$ws7.Range("I141").Value = 0
$ws7.Range("J141").Value = 1
$ws7.Range("K141").Value = 0
$ws7.Range("L141").Value = 0
$ws7.Range("M141").Value = "right"
$ws7.Range("A142").Value = 0
$ws7.Range("B142").Value = 0
$ws7.Range("C142").Value = 0.17433709262962199
$ws7.Range("D142").Value = 1
$ws7.Range("E142").Value = 1
$ws7.Range("F142").Value = -0.084066444941509194
$ws7.Range("G142").Value = 1
$ws7.Range("H142").Value = 0
$ws7.Range("I142").Value = 1
$ws7.Range("J142").Value = 0
$ws7.Range("K142").Value = 1
$ws7.Range("L142").Value = 1
$ws7.Range("M142").Value = "left"
$ws7.Range("A143").Value = 1
$ws7.Range("B143").Value = 0
$ws7.Range("C143").Value = 0.15558560892728601
$ws7.Range("D143").Value = 0
$ws7.Range("E143").Value = 0
$ws7.Range("F143").Value = -0.052135169105877503
$ws7.Range("G143").Value = 0
$ws7.Range("H143").Value = 1
$ws7.Range("I143").Value = 1
$ws7.Range("J143").Value = 1
$ws7.Range("K143").Value = 0
$ws7.Range("L143").Value = 0
$ws7.Range("M143").Value = "right"
$ws7.Range("A144").Value = 0
$ws7.Range("B144").Value = 0
$ws7.Range("C144").Value = 0.13966176000059
$ws7.Range("D144").Value = 1
$ws7.Range("E144").Value = 1
$ws7.Range("F144").Value = -0.058320479464493202
$ws7.Range("G144").Value = 1
$ws7.Range("H144").Value = 0
$ws7.Range("I144").Value = 1
$ws7.Range("J144").Value = 0
$ws7.Range("K144").Value = 1
$ws7.Range("L144").Value = 1
$ws7.Range("M144").Value = "left"
$ws7.Range("A145").Value = 1
$ws7.Range("B145").Value = 1
$ws7.Range("C145").Value = -0.062948081545613696
$ws7.Range("D145").Value = 0
$ws7.Range("E145").Value = 0
$ws7.Range("F145").Value = 0.062765628564072806
$ws7.Range("G145").Value = 0
$ws7.Range("H145").Value = 1
$ws7.Range("I145").Value = 0
$ws7.Range("J145").Value = 1
$ws7.Range("K145").Value = 0
$ws7.Range("L145").Value = 0
$ws7.Range("M145").Value = "right"
$ws7.Range("A146").Value = 0
$ws7.Range("B146").Value = 0
$ws7.Range("C146").Value = 0.076136379406990903
$ws7.Range("D146").Value = 1
$ws7.Range("E146").Value = 0
$ws7.Range("F146").Value = -0.18402015867107799
$ws7.Range("G146").Value = 1
$ws7.Range("H146").Value = 0
$ws7.Range("I146").Value = 1
$ws7.Range("J146").Value = 1
$ws7.Range("K146").Value = 1
$ws7.Range("L146").Value = 1
$ws7.Range("M146").Value = "left"
$ws7.Range("A147").Value = 0
$ws7.Range("B147").Value = 1
$ws7.Range("C147").Value = -0.238421736776155
$ws7.Range("D147").Value = 0
$ws7.Range("E147").Value = 0
$ws7.Range("F147").Value = 0.10910741960003199
$ws7.Range("G147").Value = 1
$ws7.Range("H147").Value = 1
$ws7.Range("I147").Value = 0
$ws7.Range("J147").Value = 1
$ws7.Range("K147").Value = 0
$ws7.Range("L147").Value = 0
$ws7.Range("M147").Value = "right"
$ws7.Range("A148").Value = 0
$ws7.Range("B148").Value = 0
$ws7.Range("C148").Value = 0.033508016505704799
$ws7.Range("D148").Value = 1
$ws7.Range("E148").Value = 1
$ws7.Range("F148").Value = -0.0843721776385392
$ws7.Range("G148").Value = 1
$ws7.Range("H148").Value = 0
$ws7.Range("I148").Value = 1
$ws7.Range("J148").Value = 0
$ws7.Range("K148").Value = 1
$ws7.Range("L148").Value = 1
$ws7.Range("M148").Value = "left"
$ws7.Range("A149").Value = 0
$ws7.Range("B149").Value = 1
$ws7.Range("C149").Value = -0.23516440047833001
$ws7.Range("D149").Value = 0
$ws7.Range("E149").Value = 0
$ws7.Range("F149").Value = 0.092530885586807596
$ws7.Range("G149").Value = 1
$ws7.Range("H149").Value = 1
$ws7.Range("I149").Value = 0
$ws7.Range("J149").Value = 1
$ws7.Range("K149").Value = 0
$ws7.Range("L149").Value = 0
$ws7.Range("M149").Value = "right"
$ws7.Range("A150").Value = 0
$ws7.Range("B150").Value = 0
$ws7.Range("C150").Value = 0.101548157707196
$ws7.Range("D150").Value = 1
$ws7.Range("E150").Value = 1
$ws7.Range("F150").Value = -0.107914633432386
$ws7.Range("G150").Value = 1
$ws7.Range("H150").Value = 0
$ws7.Range("I150").Value = 1
$ws7.Range("J150").Value = 0
$ws7.Range("K150").Value = 1
$ws7.Range("L150").Value = 1
$ws7.Range("M150").Value = "left"
$ws7.Range("A151").Value = 1
$ws7.Range("B151").Value = 1
$ws7.Range("C151").Value = -0.033522390965087898
$ws7.Range("D151").Value = 0
$ws7.Range("E151").Value = 0
$ws7.Range("F151").Value = 0.12868435928835401
$ws7.Range("G151").Value = 0
$ws7.Range("H151").Value = 1
$ws7.Range("I151").Value = 0
$ws7.Range("J151").Value = 1
$ws7.Range("K151").Value = 0
$ws7.Range("L151").Value = 0
$ws7.Range("M151").Value = "right"
$ws7.Range("A152").Value = 0
$ws7.Range("B152").Value = 0
$ws7.Range("C152").Value = 0.00290110895196236
$ws7.Range("D152").Value = 1
$ws7.Range("E152").Value = 1
$ws7.Range("F152").Value = -0.105256093171753
$ws7.Range("G152").Value = 1
$ws7.Range("H152").Value = 0
$ws7.Range("I152").Value = 1
$ws7.Range("J152").Value = 0
$ws7.Range("K152").Value = 1
$ws7.Range("L152").Value = 1
$ws7.Range("M152").Value = "left"
$ws7.Range("A153").Value = 0
$ws7.Range("B153").Value = 1
$ws7.Range("C153").Value = -0.25344565376419098
$ws7.Range("D153").Value = 0
$ws7.Range("E153").Value = 0
$ws7.Range("F153").Value = 0.091724783296486104
$ws7.Range("G153").Value = 1
$ws7.Range("H153").Value = 1
$ws7.Range("I153").Value = 0
$ws7.Range("J153").Value = 1
$ws7.Range("K153").Value = 0
$ws7.Range("L153").Value = 0
$ws7.Range("M153").Value = "right"
$ws7.Range("A154").Value = 0
$ws7.Range("B154").Value = 0
$ws7.Range("C154").Value = 0.074495432813834098
$ws7.Range("D154").Value = 1
$ws7.Range("E154").Value = 1
$ws7.Range("F154").Value = -0.0906559894852417
$ws7.Range("G154").Value = 1
$ws7.Range("H154").Value = 0
$ws7.Range("I154").Value = 1
$ws7.Range("J154").Value = 0
$ws7.Range("K154").Value = 1
$ws7.Range("L154").Value = 1
$ws7.Range("M154").Value = "left"
$ws7.Range("A155").Value = 0
$ws7.Range("B155").Value = 1
$ws7.Range("C155").Value = -0.25189166772856703
$ws7.Range("D155").Value = 1
$ws7.Range("E155").Value = 0
$ws7.Range("F155").Value = -0.11693544356935499
$ws7.Range("G155").Value = 1
$ws7.Range("H155").Value = 0
$ws7.Range("I155").Value = 0
$ws7.Range("J155").Value = 1
$ws7.Range("K155").Value = 0
$ws7.Range("L155").Value = 0
$ws7.Range("M155").Value = "right"
$ws7.Range("A156").Value = 0
$ws7.Range("B156").Value = 0
$ws7.Range("C156").Value = 0.0647154504782054
$ws7.Range("D156").Value = 1
$ws7.Range("E156").Value = 1
$ws7.Range("F156").Value = -0.091701864558186297
$ws7.Range("G156").Value = 1
$ws7.Range("H156").Value = 0
$ws7.Range("I156").Value = 1
$ws7.Range("J156").Value = 0
$ws7.Range("K156").Value = 1
$ws7.Range("L156").Value = 1
$ws7.Range("M156").Value = "left"
$ws7.Range("A157").Value = 1
$ws7.Range("B157").Value = 1
$ws7.Range("C157").Value = -0.055094762879053502
$ws7.Range("D157").Value = 0
$ws7.Range("E157").Value = 0
$ws7.Range("F157").Value = 0.054971607092529899
$ws7.Range("G157").Value = 0
$ws7.Range("H157").Value = 1
$ws7.Range("I157").Value = 0
$ws7.Range("J157").Value = 1
$ws7.Range("K157").Value = 0
$ws7.Range("L157").Value = 0
$ws7.Range("M157").Value = "right"

# --- New rows 131-157 on sys2_test (sheet8) ---
$ws8.Range("A131").Value = -0.023021396615663402
$ws8.Range("B131").Value = -0.14308251083017601
$ws8.Range("C131").Value = 0
$ws8.Range("D131").Value = 0
$ws8.Range("E131").Value = 0
$ws8.Range("F131").Value = "right"
$ws8.Range("A132").Value = 0.060738217814308201
$ws8.Range("B132").Value = -0.104055189098392
$ws8.Range("C132").Value = 1
$ws8.Range("D132").Value = 1
$ws8.Range("E132").Value = 1
$ws8.Range("F132").Value = "left"
$ws8.Range("A133").Value = -0.22789823849746699
$ws8.Range("B133").Value = -0.0085016880706790705
$ws8.Range("C133").Value = 1
$ws8.Range("D133").Value = 0
$ws8.Range("E133").Value = 0
$ws8.Range("F133").Value = "right"
$ws8.Range("A134").Value = -0.038013183504730999
$ws8.Range("B134").Value = -0.077462887875195702
$ws8.Range("C134").Value = 1
$ws8.Range("D134").Value = 1
$ws8.Range("E134").Value = 1
$ws8.Range("F134").Value = "left"
$ws8.Range("A135").Value = -0.053364465595300802
$ws8.Range("B135").Value = -0.0051054128281936796
$ws8.Range("C135").Value = 0
$ws8.Range("D135").Value = 0
$ws8.Range("E135").Value = 0
$ws8.Range("F135").Value = "right"
$ws8.Range("A136").Value = 0.060156846920413502
$ws8.Range("B136").Value = -0.097380944572685701
$ws8.Range("C136").Value = 1
$ws8.Range("D136").Value = 1
$ws8.Range("E136").Value = 1
$ws8.Range("F136").Value = "left"
$ws8.Range("A137").Value = -0.21314620281322399
$ws8.Range("B137").Value = 0.050115291735105302
$ws8.Range("C137").Value = 1
$ws8.Range("D137").Value = 0
$ws8.Range("E137").Value = 0
$ws8.Range("F137").Value = "right"
$ws8.Range("A138").Value = 0.091723723851644304
$ws8.Range("B138").Value = -0.072693457748126605
$ws8.Range("C138").Value = 1
$ws8.Range("D138").Value = 1
$ws8.Range("E138").Value = 1
$ws8.Range("F138").Value = "left"
$ws8.Range("A139").Value = -0.032169778198787198
$ws8.Range("B139").Value = -0.0027072716051931298
$ws8.Range("C139").Value = 0
$ws8.Range("D139").Value = 0
$ws8.Range("E139").Value = 0
$ws8.Range("F139").Value = "right"
$ws8.Range("A140").Value = 0.061722051235653497
$ws8.Range("B140").Value = -0.10005958715931799
$ws8.Range("C140").Value = 1
$ws8.Range("D140").Value = 1
$ws8.Range("E140").Value = 1
$ws8.Range("F140").Value = "left"
$ws8.Range("A141").Value = -0.0266094643094196
$ws8.Range("B141").Value = -0.14829029022474799
$ws8.Range("C141").Value = 0
$ws8.Range("D141").Value = 0
$ws8.Range("E141").Value = 0
$ws8.Range("F141").Value = "right"
$ws8.Range("A142").Value = 0.17433709262962199
$ws8.Range("B142").Value = -0.084066444941509194
$ws8.Range("C142").Value = 1
$ws8.Range("D142").Value = 1
$ws8.Range("E142").Value = 1
$ws8.Range("F142").Value = "left"
$ws8.Range("A143").Value = 0.15558560892728601
$ws8.Range("B143").Value = -0.052135169105877503
$ws8.Range("C143").Value = 0
$ws8.Range("D143").Value = 0
$ws8.Range("E143").Value = 0
$ws8.Range("F143").Value = "right"
$ws8.Range("A144").Value = 0.13966176000059
$ws8.Range("B144").Value = -0.058320479464493202
$ws8.Range("C144").Value = 1
$ws8.Range("D144").Value = 1
$ws8.Range("E144").Value = 1
$ws8.Range("F144").Value = "left"
$ws8.Range("A145").Value = -0.062948081545613696
$ws8.Range("B145").Value = 0.062765628564072806
$ws8.Range("C145").Value = 0
$ws8.Range("D145").Value = 0
$ws8.Range("E145").Value = 0
$ws8.Range("F145").Value = "right"
$ws8.Range("A146").Value = 0.076136379406990903
$ws8.Range("B146").Value = -0.18402015867107799
$ws8.Range("C146").Value = 1
$ws8.Range("D146").Value = 1
$ws8.Range("E146").Value = 1
$ws8.Range("F146").Value = "left"
$ws8.Range("A147").Value = -0.238421736776155
$ws8.Range("B147").Value = 0.10910741960003199
$ws8.Range("C147").Value = 1
$ws8.Range("D147").Value = 0
$ws8.Range("E147").Value = 0
$ws8.Range("F147").Value = "right"
$ws8.Range("A148").Value = 0.033508016505704799
$ws8.Range("B148").Value = -0.0843721776385392
$ws8.Range("C148").Value = 1
$ws8.Range("D148").Value = 1
$ws8.Range("E148").Value = 1
$ws8.Range("F148").Value = "left"
$ws8.Range("A149").Value = -0.23516440047833001
$ws8.Range("B149").Value = 0.092530885586807596
$ws8.Range("C149").Value = 1
$ws8.Range("D149").Value = 0
$ws8.Range("E149").Value = 0
$ws8.Range("F149").Value = "right"
$ws8.Range("A150").Value = 0.101548157707196
$ws8.Range("B150").Value = -0.107914633432386
$ws8.Range("C150").Value = 1
$ws8.Range("D150").Value = 1
$ws8.Range("E150").Value = 1
$ws8.Range("F150").Value = "left"
$ws8.Range("A151").Value = -0.033522390965087898
$ws8.Range("B151").Value = 0.12868435928835401
$ws8.Range("C151").Value = 0
$ws8.Range("D151").Value = 0
$ws8.Range("E151").Value = 0
$ws8.Range("F151").Value = "right"
$ws8.Range("A152").Value = 0.00290110895196236
$ws8.Range("B152").Value = -0.105256093171753
$ws8.Range("C152").Value = 1
$ws8.Range("D152").Value = 1
$ws8.Range("E152").Value = 1
$ws8.Range("F152").Value = "left"
$ws8.Range("A153").Value = -0.25344565376419098
$ws8.Range("B153").Value = 0.091724783296486104
$ws8.Range("C153").Value = 1
$ws8.Range("D153").Value = 0
$ws8.Range("E153").Value = 0
$ws8.Range("F153").Value = "right"
$ws8.Range("A154").Value = 0.074495432813834098
$ws8.Range("B154").Value = -0.0906559894852417
$ws8.Range("C154").Value = 1
$ws8.Range("D154").Value = 1
$ws8.Range("E154").Value = 1
$ws8.Range("F154").Value = "left"
$ws8.Range("A155").Value = -0.25189166772856703
$ws8.Range("B155").Value = -0.11693544356935499
$ws8.Range("C155").Value = 1
$ws8.Range("D155").Value = 0
$ws8.Range("E155").Value = 0
$ws8.Range("F155").Value = "right"
$ws8.Range("A156").Value = 0.0647154504782054
$ws8.Range("B156").Value = -0.091701864558186297
$ws8.Range("C156").Value = 1
$ws8.Range("D156").Value = 1
$ws8.Range("E156").Value = 1
$ws8.Range("F156").Value = "left"
$ws8.Range("A157").Value = -0.055094762879053502
$ws8.Range("B157").Value = 0.054971607092529899
$ws8.Range("C157").Value = 0
$ws8.Range("D157").Value = 0
$ws8.Range("E157").Value = 0
$ws8.Range("F157").Value = "right"

# --- View state: select relevant columns on sys2, then make sys2_test the
#     active sheet/tab with its own selection (mirrors the saved workbook
#     view after the edit). ---
$ws7.Activate()
$rng7 = $excel.Union($ws7.Range("C:C"), $ws7.Range("F:F"), $ws7.Range("G:G"), $ws7.Range("K:K"), $ws7.Range("L:L"), $ws7.Range("M:M"))
$rng7.Select()

$ws8.Activate()
$ws8.Range("F13").Select()
